# Rebuild the report: drop the RIMS-only measurement/identifier columns
# (DL_PSC, FREQUENCY, CPICH_POWER, TOTAL_POWER, MAX_POWER, FILENAME, RAC,
# DL_UARFCN, DC_SUPPORT, OAM_IP, SERVICE_IP), keep the surviving columns in
# place, and append two new trailing columns for province / region
# (TINH/TP, KHU VUC). Finish by bolding the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slide the header cells that survive past the columns being removed back
# into their new (contiguous) positions:
#   J (BSC_RNC_NAME)  -> E
#   K (BTS_NODEB_NAME)-> F
#   L (VENDOR)        -> G
#   N (CHECK_DATE)    -> H
#   O (CELL_CODE)     -> I
$ws.Cells.Item(2, 5).Value = $ws.Cells.Item(2, 10).Text
$ws.Cells.Item(2, 6).Value = $ws.Cells.Item(2, 11).Text
$ws.Cells.Item(2, 7).Value = $ws.Cells.Item(2, 12).Text
$ws.Cells.Item(2, 8).Value = $ws.Cells.Item(2, 14).Text
$ws.Cells.Item(2, 9).Value = $ws.Cells.Item(2, 15).Text

# Drop everything from the old column J onward (old source cells we just
# copied out of, plus every column that no longer belongs in the sheet).
$ws.Range("J1:T2").EntireColumn.Delete()

# Bring the two new trailing header cells up to the same look (font/fill/
# border) as the rest of the header row before filling them in.
$ws.Range("I2").Copy()
$ws.Range("J2:K2").PasteSpecial(-4122)
$ws.Cells.Item(2, 11).Value = "KHU VỰC"
$ws.Cells.Item(2, 10).Value = "TỈNH/TP"

$ws.Range("A1").Value = "Danh sách cell3G trên Inventory không có trên RIMS"

# Bold the whole header row.
$ws.Range("A2:K2").Font.Bold = $true

# New column widths for the rebuilt table.
$ws.Columns.Item(4).ColumnWidth = 17.5
$ws.Columns.Item(5).ColumnWidth = 19.666666666666668
$ws.Columns.Item(6).ColumnWidth = 20.166666666666668
$ws.Columns.Item(7).ColumnWidth = 10
$ws.Columns.Item(8).ColumnWidth = 12.666666666666666
$ws.Columns.Item(9).ColumnWidth = 16.5

$ws.Range("H9").Select() | Out-Null
